$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for columns I (I0) and J (IF), rows 2-16
$values = @{
    2  = @(7, 8)
    3  = @(5, 6)
    4  = @(8, 8)
    5  = @(3, 4)
    6  = @(7, 8)
    7  = @(11, 12)
    8  = @(8, 8)
    9  = @(8, 9)
    10 = @(8, 9)
    11 = @(5, 5)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(4, 4)
    15 = @(9, 9)
    16 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
